$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 499 }

# Column C holds the "Förändrad" (changed) date. Update every populated
# data row (2..lastRow) from 45192 (2023-09-23) to 45202 (2023-10-03),
# keeping the existing date formatting/style untouched.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45202
